$wb = $excel.ActiveWorkbook

# The "Status" column value "Ready for handoff" is shared (same shared
# string) by the Overview sheet's zh-cn/de-de status cells AND by the
# per-language sheets' Status cell for the 2927ab59... record. Generating
# the handback report flips all of them to "Handback transform failed".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"

# zh-cn sheet: record the Error Detail explaining the handback filename
# mismatch for the 2927ab59... record.
$wsZhCn.Range("K3").Value = "Handback file name: efik2yn0.uwv is different with handoff file name: 2927ab59-7ee9-4313-8727-3306fca3ce5e.542afa100939070f4727fe9c5990c69fd6cc1551.zh-cn."

# de-de sheet: same handback-mismatch detail for the de-de record.
$wsDeDe.Range("K3").Value = "Handback file name: efik2yn0.uwv is different with handoff file name: 2927ab59-7ee9-4313-8727-3306fca3ce5e.542afa100939070f4727fe9c5990c69fd6cc1551.de-de."
